$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new demo row (row 5): "dataset_4" - a dataset with a lot of variables,
# used to test the datatables scroller.
# Shared strings must be appended in this order: description (G), id (A), name (E)
$ws.Range("G5").Value = "a dataset with a lot of variables"
$ws.Range("A5").Value = "dataset_4"
$ws.Range("E5").Value = "dataset with lot of variables"
$ws.Range("F5").Value = "open_data"
$ws.Range("H5").Value = 200
$ws.Range("K5").NumberFormat = "@"
$ws.Range("L5").NumberFormat = "@"

# Widen columns E (name) and G (description) to fit the new, longer content;
# add width for the new column P (data_path) as well.
$ws.Columns.Item(5).ColumnWidth = 16 - 5/6
$ws.Columns.Item(7).ColumnWidth = 51 + 5/12
$ws.Columns.Item(16).ColumnWidth = 21 + 11/12

# Extend the table / autofilter to cover the new row
$table = $ws.ListObjects.Item(1)
$table.Resize($ws.Range("A1:P5"))

# Update the active view / selection
$ws.Range("E6").Select()
